$p = $ppt.ActivePresentation

# -------------------------------------------------------------------------
# 1) Slide 5: table tableStyleId change
#    {85F550FF-7759-4CDF-BFD9-3EBCBF35018E} -> {00AEA825-1DBB-414A-B3A6-189EAD490113}
# -------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{00AEA825-1DBB-414A-B3A6-189EAD490113}")

# -------------------------------------------------------------------------
# 2) Swap the two theme colour schemes (theme1.xml <-> theme2.xml content).
#    theme1.xml currently backs the Notes Master ("Office Theme" colours),
#    theme2.xml currently backs the Slide Master ("Integral"/"Red Violet"
#    colours). The edit swaps which colours are attached to which master,
#    i.e. the Slide Master ends up with the Office colours and the Notes
#    Master ends up with the Integral/Red Violet colours. (dk1/lt1 are
#    identical - 000000/FFFFFF - in both schemes, so only indices 3-12
#    need updating.)
# -------------------------------------------------------------------------
$slideMasterColors = $p.SlideMaster.Theme.ThemeColorScheme
$notesMasterColors = $p.NotesMaster.Theme.ThemeColorScheme

# Current Slide Master colours (Integral / Red Violet) - becomes Notes Master colours
$redViolet = @{
    3  = 0x514545
    4  = 0xDCD9D8
    5  = 0x912DE3
    6  = 0xCC30C8
    7  = 0xDCA64E
    8  = 0xE77547
    9  = 0xE17189
    10 = 0x7347D5
    11 = 0x259F6B
    12 = 0x8C8C8C
}

# Current Notes Master colours (Office) - becomes Slide Master colours
$office = @{
    3  = 0x6A5444
    4  = 0xE6E6E7
    5  = 0xD59B5B
    6  = 0x317DED
    7  = 0xA5A5A5
    8  = 0x00C0FF
    9  = 0xC47244
    10 = 0x47AD70
    11 = 0xC16305
    12 = 0x724F95
}

for ($i = 3; $i -le 12; $i++) {
    $slideMasterColors.Colors($i).RGB = $office[$i]
    $notesMasterColors.Colors($i).RGB = $redViolet[$i]
}
